# Disconnect charcoal production from heat fuel:
# remove the "connections" row that links the heat chain's inflow of
# biofuel to the charcoal chain's outflow of charcoal, shifting the
# remaining rows up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("connections")

# Row 29 (1-indexed, header is row 1) is:
#   B29=heat, C29=simple_heat, D29=inflow, E29=biofuel,
#   F29=charcoal, G29=outflow, H29=simple_charcoal, I29=charcoal
$ws.Rows.Item(29).Delete()

# Mirror the author's resulting selection/cursor position (purely cosmetic).
$ws.Activate()
[void]$ws.Range("B36").Select()
